# Insert two new data rows right before the current row 194, shifting
# the existing rows 194-305 down to 196-307.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(194).Resize(2).Insert()

# --- New row 194 ---
$ws.Cells.Item(194,1).Value = 9
$ws.Cells.Item(194,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(194,3).Value = "Metropolitana"
$ws.Cells.Item(194,4).Value = 44438
$ws.Cells.Item(194,5).Value = 13
$ws.Cells.Item(194,6).Value = 100112008
$ws.Cells.Item(194,7).Value = "Coliflor"
$ws.Cells.Item(194,8).Value = "Sin especificar"
$ws.Cells.Item(194,9).Value = "Primera"
$ws.Cells.Item(194,10).Value = 1600
$ws.Cells.Item(194,11).Value = 600
$ws.Cells.Item(194,12).Value = 650
$ws.Cells.Item(194,13).Value = 625
$ws.Cells.Item(194,14).Value = "`$/unidad"
$ws.Cells.Item(194,15).Value = "Región Metropolitana"
$ws.Cells.Item(194,16).Value = 625
$ws.Cells.Item(194,17).Value = 1
$ws.Cells.Item(194,18).Value = "Hortaliza"

# --- New row 195 ---
$ws.Cells.Item(195,1).Value = 9
$ws.Cells.Item(195,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(195,3).Value = "Metropolitana"
$ws.Cells.Item(195,4).Value = 44438
$ws.Cells.Item(195,5).Value = 13
$ws.Cells.Item(195,6).Value = 100112008
$ws.Cells.Item(195,7).Value = "Coliflor"
$ws.Cells.Item(195,8).Value = "Sin especificar"
$ws.Cells.Item(195,9).Value = "Segunda"
$ws.Cells.Item(195,10).Value = 610
$ws.Cells.Item(195,11).Value = 500
$ws.Cells.Item(195,12).Value = 550
$ws.Cells.Item(195,13).Value = 525
$ws.Cells.Item(195,14).Value = "`$/unidad"
$ws.Cells.Item(195,15).Value = "Región Metropolitana"
$ws.Cells.Item(195,16).Value = 525
$ws.Cells.Item(195,17).Value = 1
$ws.Cells.Item(195,18).Value = "Hortaliza"
